$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASPRgrantHistTable")

# Update the descriptive text cells so the "FY 2011-2016" ranges become "FY 2012-2016"
$ws.Range("A3").Value = "This table shows the grant awards and award dollars ASPR made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the ASPR page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars ASPR made for FY 2012-2016."
